# Fix a bug in the basis-function penalty-matrix figure on slide 3
# ("section 6"): the (row 3, col 3) entry of the first (4x4) matrix table
# was wrongly typed as "1" -- it should be "2" (the diagonal of a
# second-difference penalty matrix is 1,2,2,1).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Locate the 4x4 table shape ("Table 24") on this slide.
$table = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable -and $shp.Name -eq "Table 24") {
        $table = $shp.Table
        break
    }
}

# Fix the erroneous cell value: row 3, column 3.
$cell = $table.Cell(3, 3)
$cell.Shape.TextFrame.TextRange.Text = "2"
